# Turn "bj" into "Bj fjro", split across runs the way Word's spell-checker
# marks a flagged word: <w:proofErr spellStart/> .. <w:proofErr spellEnd/>
# around "Bj" (itself split into a "B" run and a "j" run), followed by a
# new run holding the appended " fjro" text. The existing _GoBack bookmark
# stays immediately after the text, as it was in the original paragraph.

$d = $word.ActiveDocument
$para = $d.Paragraphs.Item(1)
$r = $para.Range

$xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' w:rsidR='00FA0335' w:rsidRDefault='00582E6E'><w:proofErr w:type='spellStart'/><w:r><w:t>B</w:t></w:r><w:r><w:t>j</w:t></w:r><w:proofErr w:type='spellEnd'/><w:r><w:t xml:space='preserve'> fjro</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>"

$r.InsertXML($xml)
